$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 87161749
$ws.Range("A2").Value = 87369811
$ws.Range("A3").Value = 87369552
$ws.Range("A4").Value = 87369391
$ws.Range("A5").Value = 87369170
$ws.Range("A6").Value = 87369030
$ws.Range("A7").Value = 87368891
$ws.Range("A8").Value = 87368629
$ws.Range("A9").Value = 87145832
$ws.Range("A10").Value = 87162117
$ws.Range("A11").Value = 87364267
$ws.Range("A12").Value = 87145620
$ws.Range("A13").Value = 87145506
$ws.Range("A14").Value = 87145280
$ws.Range("A15").Value = 87145131
$ws.Range("A16").Value = 87144950
$ws.Range("A17").Value = 87144682
$ws.Range("A18").Value = 87144445
$ws.Range("A19").Value = 87162397
$ws.Range("A20").Value = 87162214
$ws.Range("A21").Value = 87162559

$ws.Range("A11:A21").RowHeight = 25.5

[void]$ws.Range("A22").Select()
